# Updated cryptos list on Thu May 18 19:00:49 UTC 2023 with GitHub Actions
#
# Applies the per-row cell updates (Price and Volume(1h) columns, plus two
# row swaps where the underlying coin ranking changed position) described
# by the upstream OOXML diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    2 = @{ D='26.646.31'; E='  -1.70%  ' }
    3 = @{ D='1.787.73'; E='  -1.73%  ' }
    4 = @{ E='  +0.10%  ' }
    5 = @{ D='307.84'; E='  -1.32%  ' }
    6 = @{ E='  +0.14%  ' }
    7 = @{ D='0.4541'; E='  +1.74%  ' }
    8 = @{ D='0.3687'; E='  -1.52%  ' }
    9 = @{ D='0.07201'; E='  -3.86%  ' }
    10 = @{ D='0.8518'; E='  -2.31%  ' }
    11 = @{ D='20.33'; E='  -2.75%  ' }
    12 = @{ D='1.787.44'; E='  -1.65%  ' }
    13 = @{ D='5.268'; E='  -1.52%  ' }
    14 = @{ D='0.07007'; E='  -1.39%  ' }
    15 = @{ D='6.450'; E='  -4.20%  ' }
    16 = @{ D='90.17'; E='  -4.44%  ' }
    17 = @{ D='1.003'; E='  +0.21%  ' }
    18 = @{ E='  -1.89%  ' }
    20 = @{ D='14.55'; E='  -3.05%  ' }
    21 = @{ D='26.650.48'; E='  -1.81%  ' }
    22 = @{ D='5.252'; E='  +0.44%  ' }
    23 = @{ D='10.54'; E='  -3.49%  ' }
    24 = @{ D='2.009.97'; E='  -1.51%  ' }
    25 = @{ D='1.901'; E='  -4.50%  ' }
    26 = @{ D='149.50'; E='  -1.45%  ' }
    27 = @{ B='LidoDAOToken'; C='https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'; D='2.130'; E='  -12.36%  ' }
    28 = @{ B='EthereumClassic'; C='https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'; D='17.99'; E='  -3.11%  ' }
    29 = @{ D='5.170'; E='  -3.02%  ' }
    30 = @{ D='113.53'; E='  -3.90%  ' }
    31 = @{ D='0.08802'; E='  -0.34%  ' }
    32 = @{ D='0.7495'; E='  -2.16%  ' }
    33 = @{ D='1.149'; E='  -2.47%  ' }
    34 = @{ D='4.423' }
    35 = @{ D='2.878'; E='  -0.40%  ' }
    36 = @{ E='  +0.08%  ' }
    37 = @{ D='1.106'; E='  +0.21%  ' }
    38 = @{ D='0.01937'; E='  -2.21%  ' }
    39 = @{ D='0.05202'; E='  -1.37%  ' }
    40 = @{ D='2.885'; E='  +1.18%  ' }
    41 = @{ D='7.092'; E='  -4.71%  ' }
    42 = @{ B='TheSandbox'; C='https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'; D='0.5174'; E='  -2.60%  ' }
    43 = @{ B='RenderToken'; C='https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'; D='2.308'; E='  +5.34%  ' }
    44 = @{ D='0.1638'; E='  -4.58%  ' }
    45 = @{ D='8.434'; E='  -3.11%  ' }
    46 = @{ D='0.4921'; E='  -2.57%  ' }
    47 = @{ D='10.22'; E='  -4.00%  ' }
    48 = @{ E='  +0.14%  ' }
    49 = @{ D='103.47'; E='  -2.01%  ' }
    50 = @{ E='  -3.73%  ' }
    51 = @{ E='  -1.47%  ' }
}

foreach ($rowNum in $changes.Keys) {
    $rowData = $changes[$rowNum]
    foreach ($col in $rowData.Keys) {
        $addr = "$col$rowNum"
        $val = $rowData[$col]
        $cell = $ws.Range($addr)

        # The Price column stores plain text (e.g. "26.646.31", "6.450")
        # rather than numbers, so trailing zeros / thousand-dot grouping
        # survive. When the replacement text parses as a bare number,
        # Excel's normal text->number inference would silently convert it
        # and drop exactly that formatting, so force the cell to Text first.
        if ($val -match '^[0-9]+(\.[0-9]+)?$') {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $val
    }
}
